$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PriceProposal")

for ($r = 3; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "Y"
}
